# Regenerate save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals.
#
# This updates column G ("K") values for rows 2-40 on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values for column G (K), keyed by row number.
$kValues = @{
    2  = 4
    3  = 0
    4  = 2
    5  = 0
    6  = 0
    7  = 0
    8  = 2
    9  = 1
    10 = 1
    11 = 0
    12 = 0
    13 = 1
    14 = 1
    15 = 1
    16 = 2
    17 = 0
    18 = 0
    19 = 0
    20 = 1
    21 = 1
    22 = 0
    23 = 2
    24 = 1
    25 = 1
    26 = 0
    27 = 0
    28 = 2
    29 = 0
    30 = 1
    31 = 2
    32 = 0
    33 = 0
    34 = 3
    35 = 1
    36 = 1
    37 = 1
    38 = 1
    39 = 2
    40 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
